# Apply updated evaluation metrics to the workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet "Summary": single data row (row 2) with updated metrics ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.3096085409252669
$wsSummary.Range("C2").Value = 0.06097560975609756
$wsSummary.Range("D2").Value = 0.8928571428571429
$wsSummary.Range("E2").Value = 0.1141552511415525
$wsSummary.Range("F2").Value = 0.239463601532567
$wsSummary.Range("G2").Value = 0.5855855855855856
$wsSummary.Range("H2").Value = 0.7602327447833066
$wsSummary.Range("I2").Value = 25
$wsSummary.Range("J2").Value = 385
$wsSummary.Range("K2").Value = 149
$wsSummary.Range("L2").Value = 3

# --- Sheet "Classification Report" ---
$wsReport = $wb.Worksheets.Item("Classification Report")

# Row 2 ("0")
$wsReport.Range("B2").Value = 0.9802631578947368
$wsReport.Range("C2").Value = 0.2790262172284644
$wsReport.Range("D2").Value = 0.434402332361516

# Row 3 ("1")
$wsReport.Range("B3").Value = 0.06097560975609756
$wsReport.Range("C3").Value = 0.8928571428571429
$wsReport.Range("D3").Value = 0.1141552511415525

# Row 4 ("accuracy")
$wsReport.Range("B4").Value = 0.3096085409252669
$wsReport.Range("C4").Value = 0.3096085409252669
$wsReport.Range("D4").Value = 0.3096085409252669
$wsReport.Range("E4").Value = 0.3096085409252669

# Row 5 ("macro avg")
$wsReport.Range("B5").Value = 0.5206193838254172
$wsReport.Range("C5").Value = 0.5859416800428037
$wsReport.Range("D5").Value = 0.2742787917515342

# Row 6 ("weighted avg")
$wsReport.Range("B6").Value = 0.9344623547846268
$wsReport.Range("C6").Value = 0.3096085409252669
$wsReport.Range("D6").Value = 0.4184469617669271

# --- Sheet "Confusion Matrix" ---
$wsMatrix = $wb.Worksheets.Item("Confusion Matrix")
$wsMatrix.Range("B2").Value = 149
$wsMatrix.Range("C2").Value = 385
$wsMatrix.Range("B3").Value = 3
$wsMatrix.Range("C3").Value = 25
